$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.843.14'
$ws.Range("E2").Value = '  -0.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.893.97'
$ws.Range("E3").Value = '  -0.46%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7903'
$ws.Range("E5").Value = '  -5.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.15'
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3145'
$ws.Range("E8").Value = '  -4.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.28'
$ws.Range("E9").Value = '  -5.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07230'
$ws.Range("E10").Value = '  +2.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08096'
$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.553'
$ws.Range("E12").Value = '  +5.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7640'
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.909.89'
$ws.Range("E14").Value = '  +0.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.37'
$ws.Range("E15").Value = '  -0.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.159'
$ws.Range("E16").Value = '  +4.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.860.25'
$ws.Range("E17").Value = '  -0.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.89'
$ws.Range("E18").Value = '  -1.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.71'
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007774'
$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.157.03'
$ws.Range("E21").Value = '  -0.67%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.128'
$ws.Range("E23").Value = '  +15.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1648'
$ws.Range("E25").Value = '  -6.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.378'
$ws.Range("E26").Value = '  +0.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.69'
$ws.Range("E27").Value = '  -1.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.69'
$ws.Range("E28").Value = '  -1.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.052'
$ws.Range("E29").Value = '  -2.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.399'
$ws.Range("E30").Value = '  +2.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.546'
$ws.Range("E31").Value = '  +1.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.460'
$ws.Range("E32").Value = '  +3.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.090'
$ws.Range("E33").Value = '  +0.25%  '

$ws.Range("E34").Value = '  -6.79%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.266'
$ws.Range("E35").Value = '  -0.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7386'
$ws.Range("E36").Value = '  +0.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("E38").Value = '  -3.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01923'
$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.781'
$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.148.20'
$ws.Range("E41").Value = '  +14.53%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.01'
$ws.Range("E42").Value = '  +1.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4419'
$ws.Range("E43").Value = '  -0.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.869'
$ws.Range("E44").Value = '  -1.60%  '

$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.20'
$ws.Range("E46").Value = '  +2.15%  '

$ws.Range("E47").Value = '  +0.18%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.987'
$ws.Range("E48").Value = '  +1.47%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.874'
$ws.Range("E49").Value = '  -1.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.435'
$ws.Range("E50").Value = '  -1.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.998'
$ws.Range("E51").Value = '  +10.14%  '
